$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baghdati")

# Copy column J (2014..2022 layout) formatting into column K, then set the 2023 values.
$ws.Range("J1:J6").Copy() | Out-Null
$ws.Range("K1:K6").PasteSpecial(-4122) | Out-Null

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1120
$ws.Range("K5").Value = 271
$ws.Range("K6").Value = 849
